$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Fix sorting issues (use depth stencil gubbins)" task is complete
# (depth stencil added to fix z-sorting), so remove its row from the ToDo
# list. This shifts all subsequent rows up by one.
$ws.Rows.Item(2).Delete()

# Match the resulting selection state (whole row 2 selected).
$ws.Range("A2:XFD2").Select()
